$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Stamp the "About" sheet with the date the file was opened/saved
# (serial 44307 = 4/21/2021), formatted as a short date (built-in m/d/yyyy).
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
